$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
$ws = $wb.Worksheets.Item("Overview")
$ws.Range("A2").Value = "58897dc8-6a13-40e0-8e08-ccb68919633c.md"
$ws.Range("B2").Value = "e2e\58897dc8-6a13-40e0-8e08-ccb68919633c.md"
$ws.Range("E2").Value = "In Translation"
$ws.Range("F2").Value = "In Translation"

$ws.Range("A3").Value = "89425222-9cec-4bd4-84f2-504f366ed51d.md"
$ws.Range("B3").Value = "e2e\89425222-9cec-4bd4-84f2-504f366ed51d.md"
$ws.Range("E3").Value = "In Translation"
$ws.Range("F3").Value = "In Translation"

$ws.Range("A4").Value = "4c5ca8e6-ad12-4d31-b8a1-64e84aa10a05.md"
$ws.Range("B4").Value = "e2e\4c5ca8e6-ad12-4d31-b8a1-64e84aa10a05.md"

foreach ($h in $ws.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq "`$B`$2") {
        $h.TextToDisplay = "e2e\58897dc8-6a13-40e0-8e08-ccb68919633c.md"
    } elseif ($addr -eq "`$B`$3") {
        $h.TextToDisplay = "e2e\89425222-9cec-4bd4-84f2-504f366ed51d.md"
    } elseif ($addr -eq "`$B`$4") {
        $h.TextToDisplay = "e2e\4c5ca8e6-ad12-4d31-b8a1-64e84aa10a05.md"
    }
}

# --- zh-cn and de-de sheets (same row shuffle pattern) ---
foreach ($sheetName in @("zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $suffix = $sheetName + ".xlf"

    $ws.Range("A2").Value = "58897dc8-6a13-40e0-8e08-ccb68919633c.md"
    $ws.Range("C2").Value = "In Translation"
    $ws.Range("G2").Value = "58897dc8-6a13-40e0-8e08-ccb68919633c.74751ef3355350c5b4d89ccc9ff8b84cca11ba62." + $suffix

    $ws.Range("A3").Value = "89425222-9cec-4bd4-84f2-504f366ed51d.md"
    $ws.Range("C3").Value = "In Translation"
    $ws.Range("G3").Value = "89425222-9cec-4bd4-84f2-504f366ed51d.1162a38999015460f315480783fc6e02f12f0519." + $suffix

    $ws.Range("A4").Value = "4c5ca8e6-ad12-4d31-b8a1-64e84aa10a05.md"
    $ws.Range("G4").Value = "4c5ca8e6-ad12-4d31-b8a1-64e84aa10a05.da9c5fcb66782dff95c13d057abd47a1ad5f7050." + $suffix

    foreach ($h in $ws.Hyperlinks) {
        $addr = $h.Range.Address()
        if ($addr -eq "`$A`$2") {
            $h.TextToDisplay = "58897dc8-6a13-40e0-8e08-ccb68919633c.md"
        } elseif ($addr -eq "`$A`$3") {
            $h.TextToDisplay = "89425222-9cec-4bd4-84f2-504f366ed51d.md"
        } elseif ($addr -eq "`$A`$4") {
            $h.TextToDisplay = "4c5ca8e6-ad12-4d31-b8a1-64e84aa10a05.md"
        }
    }
}
